$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New BOM line: Keystone Electronics 1042P battery holder.
# Values are entered in the same order the original author did (link first,
# which is why the shared-string table picks up the URL before the part
# number / description), then the rest of the row.
$ws.Cells.Item(9, 8).Value = "https://www.mouser.pl/ProductDetail/Keystone-Electronics/1042P?qs=g2rIOKKlpoboHyq0g1zn1A%3D%3D"
$ws.Cells.Item(9, 3).Value = "1042P"
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 4).NumberFormat = "0.00"
$ws.Cells.Item(9, 5).Value = 17.29
$ws.Cells.Item(9, 5).NumberFormat = "0.00"
$ws.Cells.Item(9, 7).Value = "MAUSER"
$ws.Cells.Item(9, 9).Value = "TAK"
$ws.Cells.Item(9, 10).Value = "Battery Holder"

# Widen column J slightly to fit the new "Battery Holder" text.
$ws.Columns.Item(10).ColumnWidth = 12.3

# Leave the selection on the newly-added row, like the author did.
$ws.Range("K9").Select()
